$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "students": the "no. of Students" counter in J3 is bumped
# (2 -> 3) and a new student row (Lucy Henckins) is added in row 4.
# -----------------------------------------------------------------
$wsStudents = $wb.Worksheets.Item("students")
$wsStudents.Cells.Item(3, 10).Value = 3

$wsStudents.Cells.Item(4, 1).Value = 3
$wsStudents.Cells.Item(4, 2).Value = "Lucy"
$wsStudents.Cells.Item(4, 3).Value = "Henckins"
$wsStudents.Cells.Item(4, 4).Value = "hjh"
$wsStudents.Cells.Item(4, 5).Value = "hjhj"
$wsStudents.Cells.Item(4, 6).Value = "hjhj"
$wsStudents.Cells.Item(4, 7).Value = "hjhjhj"
$wsStudents.Cells.Item(4, 8).Value = "hjh"

# -----------------------------------------------------------------
# Sheet "student_pswd": new student3/Lucy login row, and the row
# counter in G6 bumped (2 -> 3).
# -----------------------------------------------------------------
$wsPswd = $wb.Worksheets.Item("student_pswd")
$wsPswd.Cells.Item(4, 1).Value = 3
$wsPswd.Cells.Item(4, 2).Value = "student3"
$wsPswd.Cells.Item(4, 3).Value = "Lucy"
$wsPswd.Cells.Item(6, 7).Value = 3

# -----------------------------------------------------------------
# Sheet "student_courses": row 2 (Dave Lee) is replaced by Hasith
# Dewmina, and Lucy Henckins's pending enrollment request for Maths
# is duplicated across rows 3-6 (the buggy notification system
# referenced in the commit message re-writes once per pending
# request). L4 (No. of Students) goes 1 -> 5. Rows 7-13 are left as
# trailing blank rows (matching the sheet's new used range).
# -----------------------------------------------------------------
$wsCourses = $wb.Worksheets.Item("student_courses")

$wsCourses.Cells.Item(2, 1).Value = 1
$wsCourses.Cells.Item(2, 2).Value = "Hasith"
$wsCourses.Cells.Item(2, 3).Value = "Dewmina"

$wsCourses.Cells.Item(3, 1).Value = 3
$wsCourses.Cells.Item(3, 2).Value = "Lucy"
$wsCourses.Cells.Item(3, 3).Value = "Henckins"
$wsCourses.Cells.Item(3, 4).Value = 1
$wsCourses.Cells.Item(3, 5).Value = "Maths"

$wsCourses.Cells.Item(4, 1).Value = 3
$wsCourses.Cells.Item(4, 2).Value = "Lucy"
$wsCourses.Cells.Item(4, 3).Value = "Henckins"
$wsCourses.Cells.Item(4, 4).Value = 1
$wsCourses.Cells.Item(4, 5).Value = "Maths"
$wsCourses.Cells.Item(4, 12).Value = 5

$wsCourses.Cells.Item(5, 1).Value = 3
$wsCourses.Cells.Item(5, 2).Value = "Lucy"
$wsCourses.Cells.Item(5, 3).Value = "Henckins"
$wsCourses.Cells.Item(5, 4).Value = 1
$wsCourses.Cells.Item(5, 5).Value = "Maths"

$wsCourses.Cells.Item(6, 1).Value = 3
$wsCourses.Cells.Item(6, 2).Value = "Lucy"
$wsCourses.Cells.Item(6, 3).Value = "Henckins"
$wsCourses.Cells.Item(6, 4).Value = 1
$wsCourses.Cells.Item(6, 5).Value = "Maths"

# Extend the used range down to row 13 with trailing blank rows.
for ($r = 7; $r -le 13; $r++) {
    $wsCourses.Cells.Item($r, 1).Value = 0
    $wsCourses.Cells.Item($r, 1).ClearContents()
}

# -----------------------------------------------------------------
# Sheet "notifications": the pending "COURSE ENROLLMENT" request from
# Lucy to Matt is written into row 3, this becomes the active sheet
# (the notifications screen the commit message is about), and a
# trailing blank row 4 is left behind.
# -----------------------------------------------------------------
$wsNotif = $wb.Worksheets.Item("notifications")
$wsNotif.Cells.Item(3, 1).Value = 2
$wsNotif.Cells.Item(3, 2).Value = "COURSE ENROLLMENT"
$wsNotif.Cells.Item(3, 3).Value = "Lucy"
$wsNotif.Cells.Item(3, 4).Value = "Matt"
$wsNotif.Cells.Item(3, 5).Value = "Lucy would like to enroll in the Maths"

# Extend the used range down to row 4 (trailing blank row).
$wsNotif.Cells.Item(4, 1).Value = 0
$wsNotif.Cells.Item(4, 1).ClearContents()

$wsNotif.Activate()
